$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1015.22986
$ws.Range("I15").Value = 1015.22986
$ws.Range("K15").Value = 3045.68958
$ws.Range("M15").Value = -2876.68958
$ws.Range("H33").Value = 391.5
$ws.Range("I33").Value = 50
$ws.Range("J33").Value = 459.8
$ws.Range("K33").Value = 50
$ws.Range("L33").Value = 459.8
$ws.Range("M33").Value = 179
$ws.Range("N33").Value = -917.8
$ws.Range("H43").Value = 3029.1667
$ws.Range("J43").Value = 3948.5
$ws.Range("L43").Value = 3948.5
$ws.Range("N43").Value = -4086.5
$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 15000
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -16748
$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 45000
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -53736
$ws.Range("H111").Value = 1543
$ws.Range("I111").Value = 1252.6666
$ws.Range("K111").Value = 3757.9998
$ws.Range("M111").Value = -690.9998000000001
$ws.Range("H112").Value = 1330.3208
$ws.Range("I112").Value = 874
$ws.Range("K112").Value = 2622
$ws.Range("M112").Value = -1514
$ws.Range("H131").Value = 3262.6924
$ws.Range("I131").Value = 2083.125
$ws.Range("J131").Value = 5150
$ws.Range("K131").Value = 6249.375
$ws.Range("L131").Value = 15450
$ws.Range("M131").Value = -1209.375
$ws.Range("N131").Value = -25530
$ws.Range("H132").Value = 207529.69
$ws.Range("I132").Value = 3296.4285
$ws.Range("K132").Value = 9889.2855
$ws.Range("M132").Value = -7359.2855
$ws.Range("H133").Value = 48195.453
$ws.Range("J133").Value = 48195.453
$ws.Range("L133").Value = 48195.453
$ws.Range("N133").Value = -58315.453
$ws.Range("H137").Value = 3309.3333
$ws.Range("I137").Value = 1857.3334
$ws.Range("J137").Value = 5124.3335
$ws.Range("K137").Value = 5572.0002
$ws.Range("L137").Value = 15373.0005
$ws.Range("M137").Value = -3022.0002
$ws.Range("N137").Value = -20473.0005
$ws.Range("H138").Value = 5661.62
$ws.Range("I138").Value = 943.55554
$ws.Range("J138").Value = 7406.6577
$ws.Range("K138").Value = 2830.66662
$ws.Range("L138").Value = 22219.9731
$ws.Range("M138").Value = 2309.33338
$ws.Range("N138").Value = -32499.9731
$ws.Range("H140").Value = 47945.332
$ws.Range("I140").Value = 36333.332
$ws.Range("J140").Value = 49235.555
$ws.Range("K140").Value = 36333.332
$ws.Range("L140").Value = 49235.555
$ws.Range("M140").Value = -31153.332
$ws.Range("N140").Value = -59595.555

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5324.7915
$ws.Range("I32").Value = 4230.5737
$ws.Range("J32").Value = 11392.728
$ws.Range("K32").Value = 4230.5737
$ws.Range("L32").Value = 11392.728
$ws.Range("M32").Value = -3943.5737
$ws.Range("N32").Value = -11966.728
$ws.Range("H132").Value = 2135.8215
$ws.Range("I132").Value = 843.1429000000001
$ws.Range("J132").Value = 6013.857
$ws.Range("K132").Value = 2529.4287
$ws.Range("L132").Value = 18041.571
$ws.Range("M132").Value = 0.57129999999961
$ws.Range("N132").Value = -23101.571

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7187.7856
$ws.Range("I20").Value = 2567.4
$ws.Range("J20").Value = 12519
$ws.Range("K20").Value = 2567.4
$ws.Range("L20").Value = 12519
$ws.Range("M20").Value = -2320.4
$ws.Range("N20").Value = -13013
$ws.Range("H105").Value = 1700
$ws.Range("I105").Value = 1659.909
$ws.Range("K105").Value = 1659.909
$ws.Range("M105").Value = 87.09099999999989
$ws.Range("H134").Value = 2167.2144
$ws.Range("I134").Value = 1394.6364
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 4183.9092
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -1648.9092
$ws.Range("N134").Value = -20070

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22730400
$ws.Range("I31").Value = 1283.5714
$ws.Range("J31").Value = 62506350
$ws.Range("K31").Value = 1283.5714
$ws.Range("L31").Value = 62506350
$ws.Range("M31").Value = -988.5714
$ws.Range("N31").Value = -62506940
$ws.Range("H34").Value = 22730400
$ws.Range("I34").Value = 1283.5714
$ws.Range("J34").Value = 62506350
$ws.Range("K34").Value = 1283.5714
$ws.Range("L34").Value = 62506350
$ws.Range("M34").Value = -1081.5714
$ws.Range("N34").Value = -62506754
$ws.Range("H107").Value = 586.625
$ws.Range("I107").Value = 493.95
$ws.Range("J107").Value = 1050
$ws.Range("K107").Value = 493.95
$ws.Range("L107").Value = 1050
$ws.Range("M107").Value = 1426.05
$ws.Range("N107").Value = -4890

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J131").Value = 971.0294
$ws.Range("L131").Value = 2913.0882
$ws.Range("N131").Value = -12993.0882

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2616.9312
$ws.Range("I132").Value = 1544.0435
$ws.Range("K132").Value = 4632.1305
$ws.Range("M132").Value = -2102.1305

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2301.3333
$ws.Range("I16").Value = 2451
$ws.Range("J16").Value = 2002
$ws.Range("K16").Value = 2451
$ws.Range("L16").Value = 2002
$ws.Range("M16").Value = -2281
$ws.Range("N16").Value = -2342
$ws.Range("H22").Value = 102670.3
$ws.Range("I22").Value = 335333.66
$ws.Range("J22").Value = 2957.4285
$ws.Range("K22").Value = 335333.66
$ws.Range("L22").Value = 2957.4285
$ws.Range("M22").Value = -335038.66
$ws.Range("N22").Value = -3547.4285
$ws.Range("H23").Value = 2752.25
$ws.Range("I23").Value = 503.66666
$ws.Range("J23").Value = 9498
$ws.Range("K23").Value = 503.66666
$ws.Range("L23").Value = 9498
$ws.Range("M23").Value = -273.66666
$ws.Range("N23").Value = -9958
$ws.Range("H26").Value = 16001
$ws.Range("J26").Value = 29993
$ws.Range("L26").Value = 29993
$ws.Range("N26").Value = -30583
$ws.Range("H27").Value = 102670.3
$ws.Range("I27").Value = 335333.66
$ws.Range("J27").Value = 2957.4285
$ws.Range("K27").Value = 335333.66
$ws.Range("L27").Value = 2957.4285
$ws.Range("M27").Value = -335226.66
$ws.Range("N27").Value = -3171.4285
$ws.Range("H42").Value = 42497
$ws.Range("J42").Value = 42497
$ws.Range("L42").Value = 42497
$ws.Range("N42").Value = -43623
$ws.Range("H45").Value = 39946
$ws.Range("J45").Value = 39995
$ws.Range("L45").Value = 39995
$ws.Range("N45").Value = -40809
$ws.Range("H46").Value = 1825.381
$ws.Range("J46").Value = 1535.1111
$ws.Range("L46").Value = 1535.1111
$ws.Range("N46").Value = -1911.1111
$ws.Range("H49").Value = 42497
$ws.Range("J49").Value = 42497
$ws.Range("L49").Value = 42497
$ws.Range("N49").Value = -42791
$ws.Range("H54").Value = 35025.668
$ws.Range("J54").Value = 35025.668
$ws.Range("L54").Value = 35025.668
$ws.Range("N54").Value = -36313.668
$ws.Range("H93").Value = 4274895
$ws.Range("I93").Value = 7408229
$ws.Range("J93").Value = 2167.5454
$ws.Range("K93").Value = 7408229
$ws.Range("L93").Value = 2167.5454
$ws.Range("M93").Value = -7406981
$ws.Range("N93").Value = -4663.5454
$ws.Range("H123").Value = 27927.4
$ws.Range("J123").Value = 27927.4
$ws.Range("L123").Value = 27927.4
$ws.Range("N123").Value = -37727.4
$ws.Range("H136").Value = 2277.6
$ws.Range("I136").Value = 1211.2858
$ws.Range("J136").Value = 6542.857
$ws.Range("K136").Value = 3633.8574
$ws.Range("L136").Value = 19628.571
$ws.Range("M136").Value = -1083.8574
$ws.Range("N136").Value = -24728.571

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5727.7856
$ws.Range("I136").Value = 1498.7778
$ws.Range("J136").Value = 13340
$ws.Range("K136").Value = 4496.3334
$ws.Range("L136").Value = 40020
$ws.Range("M136").Value = -1946.3334
$ws.Range("N136").Value = -45120
